$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date in column C, same style as B1 (copy formats only, then set value)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C1").Value = "13-01-2023"

# Final ordered data (label, old-date value, new-date value)
$data = @(
    @("Allaria Acciones", 102326.47, 102454.54),
    @("Alpha Mega", 136797.84, 136194.07),
    @("Bull Market", 196652.38, 196810.9),
    @("Compass Crecimiento", 391972.41, 390243.02),
    @("Delta Acciones", 6737.66, 6792.31),
    @("Delta Internacional", 8493.879999999999, 8490.469999999999),
    @("Delta Latinoamerica", 12704.26, 12695.93),
    @("Delta Select", 228620.38, 229548.87),
    @("Delta gestion V", 106425.39, 107382.82),
    @("Fima PB Acciones", 4679.44, 4927.07),
    @("HF Acciones Argentinas", 9784.440000000001, 9075.68),
    @("Lombardi", 17991.65, 23446.42),
    @("Megainver", 43905.49, 43913.46),
    @("Quinquela Acciones", 140375.53, 140141.9),
    @("Toronto Trust Multimercado", 18817.61, 18675.13),
    @("Toronto trust Argy", 170525.29, 170523.4),
    @("avg", 99800.63, 100082.25),
    @("total", 1596810.12, 1601315.99)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
